$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values replacing the old Strike# proxy values for G2:G76
$kValues = @(
    2, 0, 0, 0, 1, 0, 1, 0, 2, 0, 0, 1, 1, 0, 1, 0, 3, 0, 2, 0, 0, 1, 1, 1, 0, 1, 2, 0, 0, 0, 0, 2, 1, 0, 0, 1, 1, 0, 1, 0, 1, 0, 1, 0, 1, 0, 2, 1, 1, 0, 1, 1, 0, 2, 3, 2, 0, 1, 1, 0, 0, 0, 0, 1, 0, 1, 1, 1, 0, 0, 1, 2, 0, 1, 1
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}

